$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5795122385025024
$ws.Range("B1").Value = 1.186717510223389
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.738371133804321
$ws.Range("E1").Value = 1.455952525138855
